$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the M column formulas: E*0.08 -> E*0.008
$ws.Range("M2").Formula = "=E2*0.008"
$ws.Range("M3:M8").Formula = "=E3*0.008"

# Update the selection to M2:M8, active cell M2
$ws.Range("M2:M8").Select()
